$wb = $excel.ActiveWorkbook

# The user was working on the "tax" sheet before switching away; its
# selection ends up on A3 once the new sheet takes over as active tab.
$taxSheet = $wb.Worksheets.Item("tax")
$taxSheet.Activate()
$taxSheet.Range("A3").Select() | Out-Null

# Add the new "audit" sheet after the last existing sheet (keeps it last
# in tab order, after "tax").
$ws = $wb.Worksheets.Add($null, $taxSheet)
$ws.Name = "audit"

# Populate the audit rule data.
$ws.Range("A1").Value = "Template"
$ws.Range("B1").Value = "Audit"
$ws.Range("A2").Value = "closing month"
$ws.Range("B2").Value = 12
$ws.Range("A3").Value = "audit month"
$ws.Range("B3").Value = 5
$ws.Range("A4").Value = "TargetItemType"
$ws.Range("B4").Value = "Retained earnings"

# Auto-fit column A to the new labels, and leave the selection on the
# last cell that was filled in.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Range("B4").Select() | Out-Null
